# repull data, push all data, mean calculation
# Update the dSF (column F) values for the rows whose source data was repulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = -4
    14 = -1
    19 = -2
    22 = 1
    26 = 6
    27 = 2
    34 = 0
    48 = -1
    52 = 4
    54 = -1
    58 = -4
    61 = -1
    63 = -3
    66 = 3
    69 = 0
    70 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
